# Replace "Spotify" with "olympic" in the Content Placeholder of slide 3
# (the intro sentence "Data Analysis with Spotify Data analysis dataset
# from the years 1920 onward...").
#
# Assigning the new text back onto the *whole* TextRange (rather than onto
# a Characters() sub-range) makes the host diff the old/new strings and
# re-use the surrounding runs' formatting for the untouched text, which is
# exactly the run split the target deck shows: the sentence is split into
# "Data Analysis with ", "olympic" and " Data analysis dataset from the ".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$oldText = $tr.Text
$newText = $oldText.Replace("Spotify", "olympic")

if ($newText -ne $oldText) {
    $tr.Text = $newText
}
